$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Inscritos" (column E) counts for a few rows in the data table
$ws.Range("E5").Value = 25
$ws.Range("E15").Value = 80
$ws.Range("E16").Value = 278
$ws.Range("E18").Value = 80
